$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

# Rows 2-25, columns B (2) through AE (31) all currently hold 0.5 and
# should be updated to 0.3 ("Share of Cost Effective Capacity Built in a
# Single Year" dropped from 50% to 30%).
$ws.Range("B2:AE25").Value = 0.3
